# Regenerate the "K" column (G) values on Sheet1 using freshly
# calculated s_vals (was previously derived from Strike#).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 4
    4  = 1
    5  = 7
    6  = 1
    7  = 3
    8  = 3
    9  = 8
    10 = 9
    11 = 5
    12 = 4
    13 = 7
    14 = 6
    15 = 6
    16 = 4
    17 = 4
    18 = 2
    19 = 4
    20 = 4
    21 = 2
    22 = 12
    23 = 5
    24 = 8
    25 = 6
    26 = 7
    27 = 8
    28 = 2
    29 = 2
    30 = 2
    31 = 4
    32 = 4
    33 = 9
    34 = 7
    35 = 2
    36 = 4
    37 = 1
    38 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
